# The workbook tracks one weekly price record per row for "Acelga" at
# "Vega Modelo de Temuco". A new week's record (Fecha 2021-10-15) needs to
# be inserted chronologically right after row 105 (Fecha 2021-01-13),
# pushing the existing row 106 (and everything below it) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 106; rows 106..203 shift down to 107..204.
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(106, 1).Value  = 10
$ws.Cells.Item(106, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(106, 3).Value  = "La Araucanía"
$ws.Cells.Item(106, 4).Value  = 44484
$ws.Cells.Item(106, 5).Value  = 9
$ws.Cells.Item(106, 6).Value  = 100112009
$ws.Cells.Item(106, 7).Value  = "Acelga"
$ws.Cells.Item(106, 8).Value  = "Sin especificar"
$ws.Cells.Item(106, 9).Value  = "Primera"
$ws.Cells.Item(106, 10).Value = 40
$ws.Cells.Item(106, 11).Value = 7000
$ws.Cells.Item(106, 12).Value = 8000
$ws.Cells.Item(106, 13).Value = 7500
$ws.Cells.Item(106, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(106, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(106, 16).Value = 625
$ws.Cells.Item(106, 17).Value = 12
$ws.Cells.Item(106, 18).Value = "Hortaliza"
